$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("C17").Value = 18
$ws.Range("G17").Value = 16

# Row 20
$ws.Range("D20").Value = "88.9% der Karten"

# Rows 27-31 (Karten ohne Aktivität)
$ws.Range("C27").Value = 254
$ws.Range("C28").Value = 187
$ws.Range("C29").Value = 182
$ws.Range("C30").Value = 155
$ws.Range("B31").Value = "25.05. Recruiting Abend / Team Abend"
$ws.Range("C31").Value = 119

# Row 34 (Karten ohne Mitglied)
$ws.Range("F34").Value = 4
$ws.Range("G34").Value = "(22.2%)"

# Aktivste Mitglieder
$ws.Range("C35").Value = 7
$ws.Range("C37").Value = 3
$ws.Range("C38").Value = 3
$ws.Range("C39").Value = 3
$ws.Range("F39").Value = 7
$ws.Range("F41").Value = 3
$ws.Range("F42").Value = 3
$ws.Range("F43").Value = 3

# Meiste Abgeschlossene Karten - names
$ws.Range("B48").Value = "Günther Kirchen"
$ws.Range("B50").Value = "Oliver Großheim"
$ws.Range("B51").Value = "Vanessa Raskob"
$ws.Range("B52").Value = "Matthias Bausch"

# Detailliert - Alle Karten / Sterne
$ws.Range("B70").Value = 18
$ws.Range("F70").Value = "Theresa Schmid"
$ws.Range("G70").Value = 2

$ws.Range("B71").Value = 4
$ws.Range("F71").Value = "Johanna Bracke"
$ws.Range("G71").Value = 2

$ws.Range("F72").Value = "Peter Augustin"
$ws.Range("G72").Value = 2

$ws.Range("F73").Value = "Svea Reimann"

$ws.Range("F74").Value = "Christoph Netsch"
$ws.Range("G74").Value = 1
